# Actualización automática 2025-07-25 09:55:09
# Update figures on the "CUMPLIMIENTO MENSUAL" sheet to reflect the new
# sales totals for the "OTROS" group and the grand TOTAL row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 ("OTROS" group): VENTA and POR CUMPLIR updated
$ws.Range("D2").Value = 5874.91
$ws.Range("E2").Value = -5874.91

# Row 4 (TOTAL row): VENTA, POR CUMPLIR and CUMPLIMIENTO updated
$ws.Range("D4").Value = 19843.18
$ws.Range("E4").Value = -6119.84
$ws.Range("F4").Value = 1.445943917442838
